# Assignment 2 changes: processing file, processeddata, and eda files
# - Data sheet: row 4's Height value was mis-entered as the text "sixty";
#   fix it to the numeric value 60.
# - Data sheet: row 12 was missing a Weight value; fill it in as 0.
# - Leave the workbook with the "Data" sheet active/selected (cell L11),
#   matching where the author ended up after editing.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")

# Fix the bad text entry "sixty" -> numeric 60
$ws1.Range("A4").Value = 60

# Fill in the previously empty Weight cell for row 12
$ws1.Range("B12").Value = 0

# Make "Data" the active sheet/tab with the author's final selection
$ws1.Activate() | Out-Null
$ws1.Range("L11").Select() | Out-Null
